# Reproduces the "ifForOnSameLine" fixture edit:
#  - Para 1 "ABCDEFG" -> "Scenario 1"
#  - the blank paragraph right after it is removed
#  - the big "+++IF list+++ ... +++END-IF+++" paragraph is split in two
#    (right after "INS $data+++  ", before "+++END-FOR"), and the stray
#    "_GoBack" bookmark around "list" is dropped
#  - the trailing "abcd" paragraph is replaced with a new "Scenario 2"
#    block demonstrating nested IFs on the same line, followed by a
#    trailing blank paragraph
#
# We work from the bottom of the document upward so that earlier
# paragraph indices stay valid while later ones are rewritten.

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Step 1: last paragraph ("abcd") -> "Scenario 2" block -----------------
$p5 = $d.Paragraphs.Item(5)
$body5 = '<w:body>' + `
  '<w:p><w:r><w:t>Scenario 2</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:t xml:space="preserve">+++IF list+++ +++IF </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>list[</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>0]+++ +++END-IF+++</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:t>+++END-IF+++</w:t></w:r></w:p>' + `
  '<w:p/>' + `
  '</w:body>'
[void]$p5.Range.InsertXML($pkgOpen + $body5 + $pkgClose)

# --- Step 2: split the big IF/FOR paragraph (paragraph 3) in two, ----------
#             dropping the "_GoBack" bookmark around "list"
$p3 = $d.Paragraphs.Item(3)
$body3 = '<w:body>' + `
  '<w:p>' + `
    '<w:r><w:t>+++</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">IF </w:t></w:r>' + `
    '<w:r><w:t>list</w:t></w:r>' + `
    '<w:r><w:t>+++</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>+++</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">FOR data IN </w:t></w:r>' + `
    '<w:r><w:t>list</w:t></w:r>' + `
    '<w:r><w:t>+++ +++</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">INS </w:t></w:r>' + `
    '<w:r><w:t>$</w:t></w:r>' + `
    '<w:r><w:t>data</w:t></w:r>' + `
    '<w:r><w:t>+++</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
  '</w:p>' + `
  '<w:p>' + `
    '<w:r><w:t>+++</w:t></w:r>' + `
    '<w:r><w:t>END-FOR</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> data</w:t></w:r>' + `
    '<w:r><w:t>+++ +++</w:t></w:r>' + `
    '<w:r><w:t>END-IF</w:t></w:r>' + `
    '<w:r><w:t>+++</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body>'
[void]$p3.Range.InsertXML($pkgOpen + $body3 + $pkgClose)

# --- Step 3: remove the blank paragraph right after "ABCDEFG" -------------
[void]$d.Paragraphs.Item(2).Range.Delete()

# --- Step 4: rename the first paragraph ------------------------------------
$d.Paragraphs.Item(1).Range.Text = "Scenario 1"
